$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 7).Value = 7.6715765
$ws.Cells.Item(2, 8).Value = 15.343153
$ws.Cells.Item(2, 9).Value = 0.2438217746866768
$ws.Cells.Item(2, 10).Value = 0.186760439442944
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 13).Value = 2.5411105
$ws.Cells.Item(2, 14).Value = 5.082221
$ws.Cells.Item(2, 15).Value = 0.1395356071449997
$ws.Cells.Item(2, 16).Value = 0.1016036911443679
$ws.Cells.Item(2, 17).Value = 19.49432359570325
$ws.Cells.Item(2, 18).Value = 77.97729438281299
$ws.Cells.Item(2, 19).Value = 0.03402181936607675
$ws.Cells.Item(2, 20).Value = 0.01897555000714731

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 7).Value = 7.6715765
$ws.Cells.Item(3, 8).Value = 15.343153
$ws.Cells.Item(3, 9).Value = 0.2438217746866768
$ws.Cells.Item(3, 10).Value = 0.186760439442944
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 3.190907666666666
$ws.Cells.Item(3, 14).Value = 9.572723
$ws.Cells.Item(3, 15).Value = 0.1752167954175812
$ws.Cells.Item(3, 16).Value = 0.1913777443174129
$ws.Cells.Item(3, 17).Value = 24.47929226926983
$ws.Cells.Item(3, 18).Value = 146.875753615619
$ws.Cells.Item(3, 19).Value = 0.04272167001362702
$ws.Cells.Item(3, 20).Value = 0.03574179162831943

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 7).Value = 7.6715765
$ws.Cells.Item(4, 8).Value = 15.343153
$ws.Cells.Item(4, 9).Value = 0.2438217746866768
$ws.Cells.Item(4, 10).Value = 0.186760439442944
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 4.948958666666666
$ws.Cells.Item(4, 14).Value = 14.846876
$ws.Cells.Item(4, 15).Value = 0.2717536101987069
$ws.Cells.Item(4, 16).Value = 0.2968185373211294
$ws.Cells.Item(4, 17).Value = 37.96631500667133
$ws.Cells.Item(4, 18).Value = 227.797890040028
$ws.Cells.Item(4, 19).Value = 0.0662594475161601
$ws.Cells.Item(4, 20).Value = 0.05543396046490602

$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 7).Value = 7.6715765
$ws.Cells.Item(5, 8).Value = 15.343153
$ws.Cells.Item(5, 9).Value = 0.2438217746866768
$ws.Cells.Item(5, 10).Value = 0.186760439442944
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 3.958308666666667
$ws.Cells.Item(5, 14).Value = 11.874926
$ws.Cells.Item(5, 15).Value = 0.2173557596454965
$ws.Cells.Item(5, 16).Value = 0.2374033544913186
$ws.Cells.Item(5, 17).Value = 30.36646774694633
$ws.Cells.Item(5, 18).Value = 182.198806481678
$ws.Cells.Item(5, 19).Value = 0.05299606705513572
$ws.Cells.Item(5, 20).Value = 0.04433755481002769

$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 7).Value = 7.6715765
$ws.Cells.Item(6, 8).Value = 15.343153
$ws.Cells.Item(6, 9).Value = 0.2438217746866768
$ws.Cells.Item(6, 10).Value = 0.186760439442944
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 1.499473
$ws.Cells.Item(6, 14).Value = 4.498419
$ws.Cells.Item(6, 15).Value = 0.08233796816491612
$ws.Cells.Item(6, 16).Value = 0.08993232972630592
$ws.Cells.Item(6, 17).Value = 11.5033218291845
$ws.Cells.Item(6, 18).Value = 69.019930975107
$ws.Cells.Item(6, 19).Value = 0.02007578952206494
$ws.Cells.Item(6, 20).Value = 0.01679580141981264

$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 7).Value = 7.6715765
$ws.Cells.Item(7, 8).Value = 15.343153
$ws.Cells.Item(7, 9).Value = 0.2438217746866768
$ws.Cells.Item(7, 10).Value = 0.186760439442944
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 13).Value = 2.072439
$ws.Cells.Item(7, 14).Value = 4.144878
$ws.Cells.Item(7, 15).Value = 0.1138002594282995
$ws.Cells.Item(7, 16).Value = 0.08286434299946525
$ws.Cells.Item(7, 17).Value = 15.8988743300835
$ws.Cells.Item(7, 18).Value = 63.595497320334
$ws.Cells.Item(7, 19).Value = 0.02774698121361221
$ws.Cells.Item(7, 20).Value = 0.01547578111273097

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 3.057640333333333
$ws.Cells.Item(8, 8).Value = 9.172920999999999
$ws.Cells.Item(8, 9).Value = 0.09717941187536815
$ws.Cells.Item(8, 10).Value = 0.1116549353927064
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 13).Value = 2.5411105
$ws.Cells.Item(8, 14).Value = 5.082221
$ws.Cells.Item(8, 15).Value = 0.1395356071449997
$ws.Cells.Item(8, 16).Value = 0.1016036911443679
$ws.Cells.Item(8, 17).Value = 7.769801956256832
$ws.Cells.Item(8, 18).Value = 46.61881173754099
$ws.Cells.Item(8, 19).Value = 0.01355998823802348
$ws.Cells.Item(8, 20).Value = 0.01134455357038489

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 3.057640333333333
$ws.Cells.Item(9, 8).Value = 9.172920999999999
$ws.Cells.Item(9, 9).Value = 0.09717941187536815
$ws.Cells.Item(9, 10).Value = 0.1116549353927064
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 3.190907666666666
$ws.Cells.Item(9, 14).Value = 9.572723
$ws.Cells.Item(9, 15).Value = 0.1752167954175812
$ws.Cells.Item(9, 16).Value = 0.1913777443174129
$ws.Cells.Item(9, 17).Value = 9.756647981542555
$ws.Cells.Item(9, 18).Value = 87.80983183388298
$ws.Cells.Item(9, 19).Value = 0.01702746512936724
$ws.Cells.Item(9, 20).Value = 0.02136826967736263

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 3.057640333333333
$ws.Cells.Item(10, 8).Value = 9.172920999999999
$ws.Cells.Item(10, 9).Value = 0.09717941187536815
$ws.Cells.Item(10, 10).Value = 0.1116549353927064
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 4.948958666666666
$ws.Cells.Item(10, 14).Value = 14.846876
$ws.Cells.Item(10, 15).Value = 0.2717536101987069
$ws.Cells.Item(10, 16).Value = 0.2968185373211294
$ws.Cells.Item(10, 17).Value = 15.13213562719955
$ws.Cells.Item(10, 18).Value = 136.189220644796
$ws.Cells.Item(10, 19).Value = 0.02640885601411838
$ws.Cells.Item(10, 20).Value = 0.03314125460794833

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 3.057640333333333
$ws.Cells.Item(11, 8).Value = 9.172920999999999
$ws.Cells.Item(11, 9).Value = 0.09717941187536815
$ws.Cells.Item(11, 10).Value = 0.1116549353927064
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 3.958308666666667
$ws.Cells.Item(11, 14).Value = 11.874926
$ws.Cells.Item(11, 15).Value = 0.2173557596454965
$ws.Cells.Item(11, 16).Value = 0.2374033544913186
$ws.Cells.Item(11, 17).Value = 12.10308423098289
$ws.Cells.Item(11, 18).Value = 108.927758078846
$ws.Cells.Item(11, 19).Value = 0.02112250489007323
$ws.Cells.Item(11, 20).Value = 0.02650725620773996

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 3.057640333333333
$ws.Cells.Item(12, 8).Value = 9.172920999999999
$ws.Cells.Item(12, 9).Value = 0.09717941187536815
$ws.Cells.Item(12, 10).Value = 0.1116549353927064
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 1.499473
$ws.Cells.Item(12, 14).Value = 4.498419
$ws.Cells.Item(12, 15).Value = 0.08233796816491612
$ws.Cells.Item(12, 16).Value = 0.08993232972630592
$ws.Cells.Item(12, 17).Value = 4.584849123544333
$ws.Cells.Item(12, 18).Value = 41.263642111899
$ws.Cells.Item(12, 19).Value = 0.008001555321279333
$ws.Cells.Item(12, 20).Value = 0.01004138846530626

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 3.057640333333333
$ws.Cells.Item(13, 8).Value = 9.172920999999999
$ws.Cells.Item(13, 9).Value = 0.09717941187536815
$ws.Cells.Item(13, 10).Value = 0.1116549353927064
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 13).Value = 2.072439
$ws.Cells.Item(13, 14).Value = 4.144878
$ws.Cells.Item(13, 15).Value = 0.1138002594282995
$ws.Cells.Item(13, 16).Value = 0.08286434299946525
$ws.Cells.Item(13, 17).Value = 6.336773074773
$ws.Cells.Item(13, 18).Value = 38.020638448638
$ws.Cells.Item(13, 19).Value = 0.01105904228250647
$ws.Cells.Item(13, 20).Value = 0.009252212863964356

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 4.249111
$ws.Cells.Item(14, 8).Value = 12.747333
$ws.Cells.Item(14, 9).Value = 0.1350473119652369
$ws.Cells.Item(14, 10).Value = 0.1551635125326289
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 13).Value = 2.5411105
$ws.Cells.Item(14, 14).Value = 5.082221
$ws.Cells.Item(14, 15).Value = 0.1395356071449997
$ws.Cells.Item(14, 16).Value = 0.1016036911443679
$ws.Cells.Item(14, 17).Value = 10.7974605777655
$ws.Cells.Item(14, 18).Value = 64.78476346659301
$ws.Cells.Item(14, 19).Value = 0.0188439086683695
$ws.Cells.Item(14, 20).Value = 0.01576518560424048

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 4.249111
$ws.Cells.Item(15, 8).Value = 12.747333
$ws.Cells.Item(15, 9).Value = 0.1350473119652369
$ws.Cells.Item(15, 10).Value = 0.1551635125326289
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 3.190907666666666
$ws.Cells.Item(15, 14).Value = 9.572723
$ws.Cells.Item(15, 15).Value = 0.1752167954175812
$ws.Cells.Item(15, 16).Value = 0.1913777443174129
$ws.Cells.Item(15, 17).Value = 13.55852086641767
$ws.Cells.Item(15, 18).Value = 122.026687797759
$ws.Cells.Item(15, 19).Value = 0.02366255723230717
$ws.Cells.Item(15, 20).Value = 0.02969484302886115

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 4.249111
$ws.Cells.Item(16, 8).Value = 12.747333
$ws.Cells.Item(16, 9).Value = 0.1350473119652369
$ws.Cells.Item(16, 10).Value = 0.1551635125326289
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 4.948958666666666
$ws.Cells.Item(16, 14).Value = 14.846876
$ws.Cells.Item(16, 15).Value = 0.2717536101987069
$ws.Cells.Item(16, 16).Value = 0.2968185373211294
$ws.Cells.Item(16, 17).Value = 21.02867470907866
$ws.Cells.Item(16, 18).Value = 189.258072381708
$ws.Cells.Item(16, 19).Value = 0.03669959457418415
$ws.Cells.Item(16, 20).Value = 0.04605540683554365

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 4.249111
$ws.Cells.Item(17, 8).Value = 12.747333
$ws.Cells.Item(17, 9).Value = 0.1350473119652369
$ws.Cells.Item(17, 10).Value = 0.1551635125326289
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 3.958308666666667
$ws.Cells.Item(17, 14).Value = 11.874926
$ws.Cells.Item(17, 15).Value = 0.2173557596454965
$ws.Cells.Item(17, 16).Value = 0.2374033544913186
$ws.Cells.Item(17, 17).Value = 16.81929289692867
$ws.Cells.Item(17, 18).Value = 151.373636072358
$ws.Cells.Item(17, 19).Value = 0.02935331108028641
$ws.Cells.Item(17, 20).Value = 0.03683633836990186

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 4.249111
$ws.Cells.Item(18, 8).Value = 12.747333
$ws.Cells.Item(18, 9).Value = 0.1350473119652369
$ws.Cells.Item(18, 10).Value = 0.1551635125326289
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 1.499473
$ws.Cells.Item(18, 14).Value = 4.498419
$ws.Cells.Item(18, 15).Value = 0.08233796816491612
$ws.Cells.Item(18, 16).Value = 0.08993232972630592
$ws.Cells.Item(18, 17).Value = 6.371427218503
$ws.Cells.Item(18, 18).Value = 57.342844966527
$ws.Cells.Item(18, 19).Value = 0.01111952127335117
$ws.Cells.Item(18, 20).Value = 0.01395421617057618

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 4.249111
$ws.Cells.Item(19, 8).Value = 12.747333
$ws.Cells.Item(19, 9).Value = 0.1350473119652369
$ws.Cells.Item(19, 10).Value = 0.1551635125326289
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 13).Value = 2.072439
$ws.Cells.Item(19, 14).Value = 4.144878
$ws.Cells.Item(19, 15).Value = 0.1138002594282995
$ws.Cells.Item(19, 16).Value = 0.08286434299946525
$ws.Cells.Item(19, 17).Value = 8.806023351729001
$ws.Cells.Item(19, 18).Value = 52.83614011037401
$ws.Cells.Item(19, 19).Value = 0.01536841913673846
$ws.Cells.Item(19, 20).Value = 0.01285752252350559

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 4.734553
$ws.Cells.Item(20, 8).Value = 14.203659
$ws.Cells.Item(20, 9).Value = 0.1504758656591809
$ws.Cells.Item(20, 10).Value = 0.1728902525144426
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 13).Value = 2.5411105
$ws.Cells.Item(20, 14).Value = 5.082221
$ws.Cells.Item(20, 15).Value = 0.1395356071449997
$ws.Cells.Item(20, 16).Value = 0.1016036911443679
$ws.Cells.Item(20, 17).Value = 12.0310223411065
$ws.Cells.Item(20, 18).Value = 72.186134046639
$ws.Cells.Item(20, 19).Value = 0.02099674127542322
$ws.Cells.Item(20, 20).Value = 0.01756628781834919

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 4.734553
$ws.Cells.Item(21, 8).Value = 14.203659
$ws.Cells.Item(21, 9).Value = 0.1504758656591809
$ws.Cells.Item(21, 10).Value = 0.1728902525144426
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 3.190907666666666
$ws.Cells.Item(21, 14).Value = 9.572723
$ws.Cells.Item(21, 15).Value = 0.1752167954175812
$ws.Cells.Item(21, 16).Value = 0.1913777443174129
$ws.Cells.Item(21, 17).Value = 15.10752146593967
$ws.Cells.Item(21, 18).Value = 135.967693193457
$ws.Cells.Item(21, 19).Value = 0.02636589896848814
$ws.Cells.Item(21, 20).Value = 0.03308734654068195

$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 4.734553
$ws.Cells.Item(22, 8).Value = 14.203659
$ws.Cells.Item(22, 9).Value = 0.1504758656591809
$ws.Cells.Item(22, 10).Value = 0.1728902525144426
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 4.948958666666666
$ws.Cells.Item(22, 14).Value = 14.846876
$ws.Cells.Item(22, 15).Value = 0.2717536101987069
$ws.Cells.Item(22, 16).Value = 0.2968185373211294
$ws.Cells.Item(22, 17).Value = 23.43110710214266
$ws.Cells.Item(22, 18).Value = 210.879963919284
$ws.Cells.Item(22, 19).Value = 0.04089235974065804
$ws.Cells.Item(22, 20).Value = 0.05131703186841757

$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 4.734553
$ws.Cells.Item(23, 8).Value = 14.203659
$ws.Cells.Item(23, 9).Value = 0.1504758656591809
$ws.Cells.Item(23, 10).Value = 0.1728902525144426
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 3.958308666666667
$ws.Cells.Item(23, 14).Value = 11.874926
$ws.Cells.Item(23, 15).Value = 0.2173557596454965
$ws.Cells.Item(23, 16).Value = 0.2374033544913186
$ws.Cells.Item(23, 17).Value = 18.74082217269267
$ws.Cells.Item(23, 18).Value = 168.667399554234
$ws.Cells.Item(23, 19).Value = 0.03270679608866495
$ws.Cells.Item(23, 20).Value = 0.0410447259057798

$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 4.734553
$ws.Cells.Item(24, 8).Value = 14.203659
$ws.Cells.Item(24, 9).Value = 0.1504758656591809
$ws.Cells.Item(24, 10).Value = 0.1728902525144426
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 1.499473
$ws.Cells.Item(24, 14).Value = 4.498419
$ws.Cells.Item(24, 15).Value = 0.08233796816491612
$ws.Cells.Item(24, 16).Value = 0.08993232972630592
$ws.Cells.Item(24, 17).Value = 7.099334390569
$ws.Cells.Item(24, 18).Value = 63.89400951512101
$ws.Cells.Item(24, 19).Value = 0.01238987703623384
$ws.Cells.Item(24, 20).Value = 0.01554842319559314

$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 4.734553
$ws.Cells.Item(25, 8).Value = 14.203659
$ws.Cells.Item(25, 9).Value = 0.1504758656591809
$ws.Cells.Item(25, 10).Value = 0.1728902525144426
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 13).Value = 2.072439
$ws.Cells.Item(25, 14).Value = 4.144878
$ws.Cells.Item(25, 15).Value = 0.1138002594282995
$ws.Cells.Item(25, 16).Value = 0.08286434299946525
$ws.Cells.Item(25, 17).Value = 9.812072284767
$ws.Cells.Item(25, 18).Value = 58.872433708602
$ws.Cells.Item(25, 19).Value = 0.01712419254971274
$ws.Cells.Item(25, 20).Value = 0.01432643718562093

$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 7.185148666666667
$ws.Cells.Item(26, 8).Value = 21.555446
$ws.Cells.Item(26, 9).Value = 0.2283618887583636
$ws.Cells.Item(26, 10).Value = 0.2623779198023151
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 2.5411105
$ws.Cells.Item(26, 14).Value = 5.082221
$ws.Cells.Item(26, 15).Value = 0.1395356071449997
$ws.Cells.Item(26, 16).Value = 0.1016036911443679
$ws.Cells.Item(26, 17).Value = 18.25825672092767
$ws.Cells.Item(26, 18).Value = 109.549540325566
$ws.Cells.Item(26, 19).Value = 0.03186461479667713
$ws.Cells.Item(26, 20).Value = 0.02665856512669615

$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 7).Value = 7.185148666666667
$ws.Cells.Item(27, 8).Value = 21.555446
$ws.Cells.Item(27, 9).Value = 0.2283618887583636
$ws.Cells.Item(27, 10).Value = 0.2623779198023151
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 13).Value = 3.190907666666666
$ws.Cells.Item(27, 14).Value = 9.572723
$ws.Cells.Item(27, 15).Value = 0.1752167954175812
$ws.Cells.Item(27, 16).Value = 0.1913777443174129
$ws.Cells.Item(27, 17).Value = 22.92714596660644
$ws.Cells.Item(27, 18).Value = 206.344313699458
$ws.Cells.Item(27, 19).Value = 0.04001283834374662
$ws.Cells.Item(27, 20).Value = 0.05021329445046214

$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 7).Value = 7.185148666666667
$ws.Cells.Item(28, 8).Value = 21.555446
$ws.Cells.Item(28, 9).Value = 0.2283618887583636
$ws.Cells.Item(28, 10).Value = 0.2623779198023151
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 13).Value = 4.948958666666666
$ws.Cells.Item(28, 14).Value = 14.846876
$ws.Cells.Item(28, 15).Value = 0.2717536101987069
$ws.Cells.Item(28, 16).Value = 0.2968185373211294
$ws.Cells.Item(28, 17).Value = 35.55900376518844
$ws.Cells.Item(28, 18).Value = 320.031033886696
$ws.Cells.Item(28, 19).Value = 0.0620581677018808
$ws.Cells.Item(28, 20).Value = 0.07787863038108378

$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 7).Value = 7.185148666666667
$ws.Cells.Item(29, 8).Value = 21.555446
$ws.Cells.Item(29, 9).Value = 0.2283618887583636
$ws.Cells.Item(29, 10).Value = 0.2623779198023151
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 13).Value = 3.958308666666667
$ws.Cells.Item(29, 14).Value = 11.874926
$ws.Cells.Item(29, 15).Value = 0.2173557596454965
$ws.Cells.Item(29, 16).Value = 0.2374033544913186
$ws.Cells.Item(29, 17).Value = 28.44103623855511
$ws.Cells.Item(29, 18).Value = 255.969326146996
$ws.Cells.Item(29, 19).Value = 0.04963577180515448
$ws.Cells.Item(29, 20).Value = 0.06228939830552379

$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 7).Value = 7.185148666666667
$ws.Cells.Item(30, 8).Value = 21.555446
$ws.Cells.Item(30, 9).Value = 0.2283618887583636
$ws.Cells.Item(30, 10).Value = 0.2623779198023151
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 13).Value = 1.499473
$ws.Cells.Item(30, 14).Value = 4.498419
$ws.Cells.Item(30, 15).Value = 0.08233796816491612
$ws.Cells.Item(30, 16).Value = 0.08993232972630592
$ws.Cells.Item(30, 17).Value = 10.77393642665267
$ws.Cells.Item(30, 18).Value = 96.965427839874
$ws.Cells.Item(30, 19).Value = 0.01880285392666626
$ws.Cells.Item(30, 20).Value = 0.02359625759656405

$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 7).Value = 7.185148666666667
$ws.Cells.Item(31, 8).Value = 21.555446
$ws.Cells.Item(31, 9).Value = 0.2283618887583636
$ws.Cells.Item(31, 10).Value = 0.2623779198023151
$ws.Cells.Item(31, 11).Value = 2
$ws.Cells.Item(31, 13).Value = 2.072439
$ws.Cells.Item(31, 14).Value = 4.144878
$ws.Cells.Item(31, 15).Value = 0.1138002594282995
$ws.Cells.Item(31, 16).Value = 0.08286434299946525
$ws.Cells.Item(31, 17).Value = 14.890782317598
$ws.Cells.Item(31, 18).Value = 89.344693905588
$ws.Cells.Item(31, 19).Value = 0.02598764218423825
$ws.Cells.Item(31, 20).Value = 0.02174177394198523

$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 7).Value = 4.56584
$ws.Cells.Item(32, 8).Value = 9.131679999999999
$ws.Cells.Item(32, 9).Value = 0.1451137470551739
$ws.Cells.Item(32, 10).Value = 0.1111529403149629
$ws.Cells.Item(32, 11).Value = 2
$ws.Cells.Item(32, 13).Value = 2.5411105
$ws.Cells.Item(32, 14).Value = 5.082221
$ws.Cells.Item(32, 15).Value = 0.1395356071449997
$ws.Cells.Item(32, 16).Value = 0.1016036911443679
$ws.Cells.Item(32, 17).Value = 11.60230396532
$ws.Cells.Item(32, 18).Value = 46.40921586128
$ws.Cells.Item(32, 19).Value = 0.02024853480042959
$ws.Cells.Item(32, 20).Value = 0.01129354901754984

$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 7).Value = 4.56584
$ws.Cells.Item(33, 8).Value = 9.131679999999999
$ws.Cells.Item(33, 9).Value = 0.1451137470551739
$ws.Cells.Item(33, 10).Value = 0.1111529403149629
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 13).Value = 3.190907666666666
$ws.Cells.Item(33, 14).Value = 9.572723
$ws.Cells.Item(33, 15).Value = 0.1752167954175812
$ws.Cells.Item(33, 16).Value = 0.1913777443174129
$ws.Cells.Item(33, 17).Value = 14.56917386077333
$ws.Cells.Item(33, 18).Value = 87.41504316464
$ws.Cells.Item(33, 19).Value = 0.02542636573004503
$ws.Cells.Item(33, 20).Value = 0.02127219899172562

$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(34, 7).Value = 4.56584
$ws.Cells.Item(34, 8).Value = 9.131679999999999
$ws.Cells.Item(34, 9).Value = 0.1451137470551739
$ws.Cells.Item(34, 10).Value = 0.1111529403149629
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 13).Value = 4.948958666666666
$ws.Cells.Item(34, 14).Value = 14.846876
$ws.Cells.Item(34, 15).Value = 0.2717536101987069
$ws.Cells.Item(34, 16).Value = 0.2968185373211294
$ws.Cells.Item(34, 17).Value = 22.59615343861333
$ws.Cells.Item(34, 18).Value = 135.57692063168
$ws.Cells.Item(34, 19).Value = 0.03943518465170548
$ws.Cells.Item(34, 20).Value = 0.03299225316323007

$ws.Cells.Item(35, 5).Value = 2
$ws.Cells.Item(35, 7).Value = 4.56584
$ws.Cells.Item(35, 8).Value = 9.131679999999999
$ws.Cells.Item(35, 9).Value = 0.1451137470551739
$ws.Cells.Item(35, 10).Value = 0.1111529403149629
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 13).Value = 3.958308666666667
$ws.Cells.Item(35, 14).Value = 11.874926
$ws.Cells.Item(35, 15).Value = 0.2173557596454965
$ws.Cells.Item(35, 16).Value = 0.2374033544913186
$ws.Cells.Item(35, 17).Value = 18.07300404261333
$ws.Cells.Item(35, 18).Value = 108.43802425568
$ws.Cells.Item(35, 19).Value = 0.03154130872618175
$ws.Cells.Item(35, 20).Value = 0.02638808089234551

$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 7).Value = 4.56584
$ws.Cells.Item(36, 8).Value = 9.131679999999999
$ws.Cells.Item(36, 9).Value = 0.1451137470551739
$ws.Cells.Item(36, 10).Value = 0.1111529403149629
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 13).Value = 1.499473
$ws.Cells.Item(36, 14).Value = 4.498419
$ws.Cells.Item(36, 15).Value = 0.08233796816491612
$ws.Cells.Item(36, 16).Value = 0.08993232972630592
$ws.Cells.Item(36, 17).Value = 6.846353802319999
$ws.Cells.Item(36, 18).Value = 41.07812281392
$ws.Cells.Item(36, 19).Value = 0.0119483710853206
$ws.Cells.Item(36, 20).Value = 0.009996242878453642

$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 7).Value = 4.56584
$ws.Cells.Item(37, 8).Value = 9.131679999999999
$ws.Cells.Item(37, 9).Value = 0.1451137470551739
$ws.Cells.Item(37, 10).Value = 0.1111529403149629
$ws.Cells.Item(37, 11).Value = 2
$ws.Cells.Item(37, 13).Value = 2.072439
$ws.Cells.Item(37, 14).Value = 4.144878
$ws.Cells.Item(37, 15).Value = 0.1138002594282995
$ws.Cells.Item(37, 16).Value = 0.08286434299946525
$ws.Cells.Item(37, 17).Value = 9.462424883760001
$ws.Cells.Item(37, 18).Value = 37.84969953504
$ws.Cells.Item(37, 19).Value = 0.01651398206149143
$ws.Cells.Item(37, 20).Value = 0.009210615371658171
